$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
# B2: name changed
$ws.Range("B2").Value = "Ikki maru"

# C2: phone column switched from a text value to a plain number
$ws.Range("C2").Value = 992907510905

# D2: email changed
$ws.Range("D2").Value = "adasdad@asdad.ru"

# E2 ("Ismat.") is unchanged

# --- Append new row 3 ---
$ws.Range("A3").Value = 974794263
$ws.Range("B3").Value = "Buzurgmehr Abdulloev"

# C3 keeps its leading "+" so it must stay text, not be coerced to a number
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "+992938636344"

$ws.Range("D3").Value = "123@gmail.com"
$ws.Range("E3").Value = "BuzurgmehrAbdulloev"
